# Auto-generated edit script applying numeric updates to Leve profit columns (H-N)
# across multiple crafting-job worksheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2489.7437
$ws.Range("I64").Value = 2450
$ws.Range("J64").Value = 2553.3333
$ws.Range("K64").Value = 2450
$ws.Range("L64").Value = 2553.3333
$ws.Range("M64").Value = -2202
$ws.Range("N64").Value = -3049.3333
$ws.Range("H67").Value = 2489.7437
$ws.Range("I67").Value = 2450
$ws.Range("J67").Value = 2553.3333
$ws.Range("K67").Value = 2450
$ws.Range("L67").Value = 2553.3333
$ws.Range("M67").Value = -1592
$ws.Range("N67").Value = -4269.3333
$ws.Range("H74").Value = 3187
$ws.Range("I74").Value = 2784.111
$ws.Range("K74").Value = 2784.111
$ws.Range("M74").Value = -1848.111
$ws.Range("H77").Value = 3187
$ws.Range("I77").Value = 2784.111
$ws.Range("K77").Value = 13920.555
$ws.Range("M77").Value = -9240.555
$ws.Range("H92").Value = 406.25
$ws.Range("I92").Value = 433.2143
$ws.Range("K92").Value = 433.2143
$ws.Range("M92").Value = 814.7857
$ws.Range("H100").Value = 6494.4736
$ws.Range("I100").Value = 9433.75
$ws.Range("J100").Value = 1455.7142
$ws.Range("K100").Value = 9433.75
$ws.Range("L100").Value = 1455.7142
$ws.Range("M100").Value = -8892.75
$ws.Range("N100").Value = -2537.7142
$ws.Range("H138").Value = 662.1892
$ws.Range("I138").Value = 662.1892
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 1986.5676
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 3153.4324
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13517537
$ws.Range("I32").Value = 4062.5166
$ws.Range("J32").Value = 71432430
$ws.Range("K32").Value = 4062.5166
$ws.Range("L32").Value = 71432430
$ws.Range("M32").Value = -3775.5166
$ws.Range("N32").Value = -71433004
$ws.Range("H102").Value = 1833.3334
$ws.Range("I102").Value = 1833.3334
$ws.Range("K102").Value = 1833.3334
$ws.Range("M102").Value = -211.3334
$ws.Range("H122").Value = 1102.6666
$ws.Range("I122").Value = 1102.6666
$ws.Range("K122").Value = 3307.9998
$ws.Range("M122").Value = -857.9998000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1792985
$ws.Range("I86").Value = 7252.5
$ws.Range("J86").Value = 2586644
$ws.Range("K86").Value = 7252.5
$ws.Range("L86").Value = 2586644
$ws.Range("M86").Value = -6129.5
$ws.Range("N86").Value = -2588890
$ws.Range("H89").Value = 1792985
$ws.Range("I89").Value = 7252.5
$ws.Range("J89").Value = 2586644
$ws.Range("K89").Value = 36262.5
$ws.Range("L89").Value = 12933220
$ws.Range("M89").Value = -30646.5
$ws.Range("N89").Value = -12944452
$ws.Range("H99").Value = 1769
$ws.Range("I99").Value = 975
$ws.Range("K99").Value = 975
$ws.Range("M99").Value = 523
$ws.Range("H105").Value = 90910520
$ws.Range("I105").Value = 1571
$ws.Range("J105").Value = 1000000000
$ws.Range("K105").Value = 1571
$ws.Range("L105").Value = 1000000000
$ws.Range("M105").Value = 176
$ws.Range("N105").Value = -1000003494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1313.8889
$ws.Range("I31").Value = 896.525
$ws.Range("J31").Value = 2039.7391
$ws.Range("K31").Value = 896.525
$ws.Range("L31").Value = 2039.7391
$ws.Range("M31").Value = -601.525
$ws.Range("N31").Value = -2629.7391
$ws.Range("H34").Value = 1313.8889
$ws.Range("I34").Value = 896.525
$ws.Range("J34").Value = 2039.7391
$ws.Range("K34").Value = 896.525
$ws.Range("L34").Value = 2039.7391
$ws.Range("M34").Value = -694.525
$ws.Range("N34").Value = -2443.7391
$ws.Range("H62").Value = 3998.6
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 4664.3335
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 4664.3335
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -5912.3335
$ws.Range("H65").Value = 3998.6
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 4664.3335
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 23321.6675
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -29561.6675
$ws.Range("H105").Value = 22708
$ws.Range("I105").Value = 35250
$ws.Range("K105").Value = 35250
$ws.Range("M105").Value = -33503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H97").Value = 727.0303
$ws.Range("I97").Value = 626
$ws.Range("J97").Value = 801.4737
$ws.Range("K97").Value = 626
$ws.Range("L97").Value = 801.4737
$ws.Range("M97").Value = -130
$ws.Range("N97").Value = -1793.4737
$ws.Range("H132").Value = 4446.1836
$ws.Range("I132").Value = 2038.75
$ws.Range("J132").Value = 11112.923
$ws.Range("K132").Value = 6116.25
$ws.Range("L132").Value = 33338.769
$ws.Range("M132").Value = -3586.25
$ws.Range("N132").Value = -38398.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 1744
$ws.Range("I30").Value = 658.6667
$ws.Range("K30").Value = 658.6667
$ws.Range("M30").Value = -550.6667
$ws.Range("H93").Value = 761.48
$ws.Range("I93").Value = 580.8946999999999
$ws.Range("J93").Value = 1333.3334
$ws.Range("K93").Value = 580.8946999999999
$ws.Range("L93").Value = 1333.3334
$ws.Range("M93").Value = 667.1053000000001
$ws.Range("N93").Value = -3829.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 43210
$ws.Range("J46").Value = 43210
$ws.Range("L46").Value = 43210
$ws.Range("N46").Value = -43672
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H62").Value = 5328.5713
$ws.Range("I62").Value = 4366.6665
$ws.Range("J62").Value = 6050
$ws.Range("K62").Value = 4366.6665
$ws.Range("L62").Value = 6050
$ws.Range("M62").Value = -3742.6665
$ws.Range("N62").Value = -7298
$ws.Range("H65").Value = 5328.5713
$ws.Range("I65").Value = 4366.6665
$ws.Range("J65").Value = 6050
$ws.Range("K65").Value = 21833.3325
$ws.Range("L65").Value = 30250
$ws.Range("M65").Value = -18713.3325
$ws.Range("N65").Value = -36490
$ws.Range("H81").Value = 1071.5555
$ws.Range("I81").Value = 977.2857
$ws.Range("J81").Value = 1401.5
$ws.Range("K81").Value = 1954.5714
$ws.Range("L81").Value = 2803
$ws.Range("M81").Value = -893.5714
$ws.Range("N81").Value = -4925
$ws.Range("H84").Value = 1071.5555
$ws.Range("I84").Value = 977.2857
$ws.Range("J84").Value = 1401.5
$ws.Range("K84").Value = 9772.857
$ws.Range("L84").Value = 14015
$ws.Range("M84").Value = -4468.857
$ws.Range("N84").Value = -24623
$ws.Range("H96").Value = 2588.7778
$ws.Range("J96").Value = 2885.5715
$ws.Range("L96").Value = 2885.5715
$ws.Range("N96").Value = -5631.5715
$ws.Range("H132").Value = 16210.185
$ws.Range("I132").Value = 17195.39
$ws.Range("J132").Value = 10955.75
$ws.Range("K132").Value = 51586.17
$ws.Range("L132").Value = 32867.25
$ws.Range("M132").Value = -49056.17
$ws.Range("N132").Value = -37927.25
$ws.Range("H134").Value = 43210
$ws.Range("J134").Value = 43210
$ws.Range("L134").Value = 129630
$ws.Range("N134").Value = -134700
$ws.Range("H136").Value = 1184.6296
$ws.Range("I136").Value = 709.7560999999999
$ws.Range("J136").Value = 2682.3076
$ws.Range("K136").Value = 2129.2683
$ws.Range("L136").Value = 8046.9228
$ws.Range("M136").Value = 420.7317000000003
$ws.Range("N136").Value = -13146.9228

Write-Output "Applied all Leve profit updates"
